$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '51.393.58'
Set-TextValue 'E2' '  -1.02%  '

Set-TextValue 'D3' '2.776.72'
Set-TextValue 'E3' '  -0.44%  '

Set-TextValue 'E4' '  +0.01%  '

Set-TextValue 'D5' '352.60'
Set-TextValue 'E5' '  -2.31%  '

Set-TextValue 'D6' '107.97'
Set-TextValue 'E6' '  -1.32%  '

Set-TextValue 'D7' '0.549'
Set-TextValue 'E7' '  -1.76%  '

Set-TextValue 'D8' '1.00'
Set-TextValue 'E8' '  +0.03%  '

Set-TextValue 'D9' '0.616'
Set-TextValue 'E9' '  +4.34%  '

Set-TextValue 'D10' '39.11'
Set-TextValue 'E10' '  -2.44%  '

Set-TextValue 'E11' '  +1.48%  '

Set-TextValue 'E12' '  -1.88%  '

Set-TextValue 'D13' '19.86'
Set-TextValue 'E13' '  +1.79%  '

Set-TextValue 'D14' '7.74'
Set-TextValue 'E14' '  +2.38%  '

Set-TextValue 'D15' '3.208.78'
Set-TextValue 'E15' '  -0.56%  '

Set-TextValue 'D16' '2.792.85'
Set-TextValue 'E16' '  +0.01%  '

Set-TextValue 'E17' '  -1.64%  '

Set-TextValue 'D18' '51.344.86'
Set-TextValue 'E18' '  -1.06%  '

Set-TextValue 'D19' '7.70'
Set-TextValue 'E19' '  +3.20%  '

Set-TextValue 'E20' '  +0.44%  '

Set-TextValue 'D21' '13.34'
Set-TextValue 'E21' '  +1.88%  '

Set-TextValue 'D22' '0.0₃0965'
Set-TextValue 'E22' '  -1.09%  '

Set-TextValue 'D23' '70.39'
Set-TextValue 'E23' '  +0.05%  '

Set-TextValue 'D24' '266.02'
Set-TextValue 'E24' '  -1.34%  '

Set-TextValue 'E25' '  -0.05%  '

Set-TextValue 'D26' '1.00'
Set-TextValue 'E26' '  -0.02%  '

Set-TextValue 'D27' '25.80'
Set-TextValue 'E27' '  -2.72%  '

Set-TextValue 'E28' '  +1.69%  '

Set-TextValue 'D29' '10.25'
Set-TextValue 'E29' '  -0.22%  '

Set-TextValue 'D30' '37.02'
Set-TextValue 'E30' '  +8.10%  '

Set-TextValue 'E31' '  -0.62%  '

Set-TextValue 'D32' '6.23'
Set-TextValue 'E32' '  +8.56%  '

Set-TextValue 'D33' '51.75'
Set-TextValue 'E33' '  -0.47%  '

Set-TextValue 'B34' 'RenderToken'
Set-TextValue 'C34' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D34' '5.64'
Set-TextValue 'E34' '  +7.73%  '

Set-TextValue 'B35' 'VeChain'
Set-TextValue 'C35' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D35' '0.0441'
Set-TextValue 'E35' '  -6.10%  '

Set-TextValue 'E36' '  -0.06%  '

Set-TextValue 'D37' '0.0838'
Set-TextValue 'E37' '  -0.89%  '

Set-TextValue 'D38' '18.43'
Set-TextValue 'E38' '  -3.18%  '

Set-TextValue 'D39' '3.11'
Set-TextValue 'E39' '  -3.18%  '

Set-TextValue 'D40' '1.96'
Set-TextValue 'E40' '  -1.76%  '

Set-TextValue 'E41' '  -1.23%  '

Set-TextValue 'E42' '  -5.04%  '

Set-TextValue 'D43' '120.03'
Set-TextValue 'E43' '  +0.77%  '

Set-TextValue 'D44' '2.19'
Set-TextValue 'E44' '  -2.74%  '

Set-TextValue 'D45' '21.86'
Set-TextValue 'E45' '  -0.57%  '

Set-TextValue 'D46' '2.130.68'
Set-TextValue 'E46' '  +2.28%  '

Set-TextValue 'D47' '3.34'
Set-TextValue 'E47' '  +2.83%  '

Set-TextValue 'E48' '  +5.69%  '

Set-TextValue 'D49' '0.226'
Set-TextValue 'E49' '  +18.46%  '

Set-TextValue 'D50' '5.45'
Set-TextValue 'E50' '  -5.81%  '

Set-TextValue 'E51' '  +8.35%  '
